# Auto-generated edit script: updates Leve profit-calculation values
# across multiple job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR),
# reflecting refreshed Universalis market-price data pulled by the
# scheduled runner. Values touched: currentAveragePrice(NQ/HQ),
# LevePrice(NQ/HQ) and the derived LeveProfit(NQ/HQ) columns.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether / Ether
$ws.Range("H15").Value = 1744.2
$ws.Range("I15").Value = 1744.2
$ws.Range("K15").Value = 5232.6
$ws.Range("M15").Value = -5063.6

# Row 17: One for the Road / Potion
$ws.Range("H17").Value = 789154.75
$ws.Range("J17").Value = 1025451.2
$ws.Range("L17").Value = 3076353.6
$ws.Range("N17").Value = -3076689.6

# Row 74: Adhesive of Antipathy / Wing Glue
$ws.Range("H74").Value = 7779.4
$ws.Range("I74").Value = 7724.25
$ws.Range("J74").Value = 8000
$ws.Range("K74").Value = 7724.25
$ws.Range("L74").Value = 8000
$ws.Range("M74").Value = -6788.25
$ws.Range("N74").Value = -9872

# Row 77: It's Gonna Grow Back (L) / Wing Glue
$ws.Range("H77").Value = 7779.4
$ws.Range("I77").Value = 7724.25
$ws.Range("J77").Value = 8000
$ws.Range("K77").Value = 38621.25
$ws.Range("L77").Value = 40000
$ws.Range("M77").Value = -33941.25
$ws.Range("N77").Value = -49360

# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 2657947.2
$ws.Range("I98").Value = 2843128.5
$ws.Range("K98").Value = 2843128.5
$ws.Range("M98").Value = -2841630.5

# Row 113: Amaro Kart / Starch Glue
$ws.Range("H113").Value = 3910
$ws.Range("J113").Value = 999
$ws.Range("L113").Value = 999
$ws.Range("N113").Value = -7507

# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 2657947.2
$ws.Range("I122").Value = 2843128.5
$ws.Range("K122").Value = 8529385.5
$ws.Range("M122").Value = -8526935.5

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 2309.7778
$ws.Range("I132").Value = 2186.5667
$ws.Range("K132").Value = 6559.7001
$ws.Range("M132").Value = -4029.7001

$ws = $wb.Worksheets.Item("ARM")
# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 2511.111
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 2511.111
$ws.Range("K45").Value = 0
$ws.Range("M45").Value = 2511.111
$ws.Range("N45").Value = -3265.111
$ws.Range("L45").ClearContents()

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 2101.4648
$ws.Range("I61").Value = 1370.8853
$ws.Range("K61").Value = 1370.8853
$ws.Range("M61").Value = -1158.8853

# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 2017.25
$ws.Range("I102").Value = 1593.4166
$ws.Range("J102").Value = 4560.25
$ws.Range("K102").Value = 1593.4166
$ws.Range("L102").Value = 4560.25
$ws.Range("M102").Value = 28.58339999999998
$ws.Range("N102").Value = -7804.25

# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 11221.1
$ws.Range("I110").Value = 11254.2
$ws.Range("K110").Value = 11254.2
$ws.Range("M110").Value = -9209.200000000001

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 5788.931
$ws.Range("I122").Value = 4711.16
$ws.Range("J122").Value = 12525
$ws.Range("K122").Value = 14133.48
$ws.Range("L122").Value = 37575
$ws.Range("M122").Value = -11683.48
$ws.Range("N122").Value = -42475

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 3712.2917
$ws.Range("I132").Value = 2632.8
$ws.Range("K132").Value = 7898.400000000001
$ws.Range("M132").Value = -5368.400000000001

# Row 135: Forgiveness for My Shins / Ruthenium Sabatons of Fending
$ws.Range("H135").Value = 32219.75
$ws.Range("J135").Value = 32219.75
$ws.Range("L135").Value = 32219.75
$ws.Range("N135").Value = -42359.75

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 2101.4648
$ws.Range("I136").Value = 1370.8853
$ws.Range("K136").Value = 4112.6559
$ws.Range("M136").Value = -1562.6559

$ws = $wb.Worksheets.Item("BSM")
# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 6919.25
$ws.Range("I99").Value = 7679.722
$ws.Range("J99").Value = 5550.4
$ws.Range("K99").Value = 7679.722
$ws.Range("L99").Value = 5550.4
$ws.Range("M99").Value = -6181.722
$ws.Range("N99").Value = -8546.4

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 1835.4305
$ws.Range("I134").Value = 1387.2325
$ws.Range("K134").Value = 4161.6975
$ws.Range("M134").Value = -1626.6975

$ws = $wb.Worksheets.Item("CRP")
# Row 47: Grippy When Wet / Mythril Cavalry Bow
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()

# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 1567.5454
$ws.Range("I122").Value = 1276.7894
$ws.Range("J122").Value = 3409
$ws.Range("K122").Value = 3830.3682
$ws.Range("L122").Value = 10227
$ws.Range("M122").Value = -1380.3682
$ws.Range("N122").Value = -15127

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 5383.8
$ws.Range("I132").Value = 5153.75
$ws.Range("J132").Value = 6304
$ws.Range("K132").Value = 15461.25
$ws.Range("L132").Value = 18912
$ws.Range("M132").Value = -12931.25
$ws.Range("N132").Value = -23972

$ws = $wb.Worksheets.Item("CUL")
# Row 137: Creative Chocolate / Gateau au Chocolat
$ws.Range("H137").Value = 15876325
$ws.Range("J137").Value = 23813442
$ws.Range("L137").Value = 71440326
$ws.Range("N137").Value = -71450526

$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws.Range("H97").Value = 646.9091
$ws.Range("I97").Value = 564.2
$ws.Range("K97").Value = 564.2
$ws.Range("M97").Value = -68.20000000000005

# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 2683.7693
$ws.Range("I102").Value = 2391.16
$ws.Range("K102").Value = 2391.16
$ws.Range("M102").Value = -769.1599999999999

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 1203.7273
$ws.Range("I122").Value = 1316.4445
$ws.Range("J122").Value = 696.5
$ws.Range("K122").Value = 3949.3335
$ws.Range("L122").Value = 2089.5
$ws.Range("M122").Value = -1499.3335
$ws.Range("N122").Value = -6989.5

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 28621.744
$ws.Range("I132").Value = 30801.883
$ws.Range("K132").Value = 92405.649
$ws.Range("M132").Value = -89875.649

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore / Hard Leather
$ws.Range("H16").Value = 3197.5386
$ws.Range("I16").Value = 1824.24
$ws.Range("K16").Value = 1824.24
$ws.Range("M16").Value = -1654.24

# Row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 474434.22
$ws.Range("I40").Value = 566111.4
$ws.Range("J40").Value = 2951.7144
$ws.Range("K40").Value = 566111.4
$ws.Range("L40").Value = 2951.7144
$ws.Range("M40").Value = -565975.4
$ws.Range("N40").Value = -3223.7144

# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 1376.6
$ws.Range("I61").Value = 1376.6
$ws.Range("K61").Value = 1376.6
$ws.Range("M61").Value = -1174.6

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 1376.6
$ws.Range("I113").Value = 1376.6
$ws.Range("K113").Value = 1376.6
$ws.Range("M113").Value = 793.4000000000001

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 3883
$ws.Range("I122").Value = 3590.5557
$ws.Range("K122").Value = 10771.6671
$ws.Range("M122").Value = -8321.667099999999

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 6289.234
$ws.Range("I132").Value = 2785.7
$ws.Range("J132").Value = 8884.444
$ws.Range("K132").Value = 8357.099999999999
$ws.Range("L132").Value = 26653.332
$ws.Range("M132").Value = -5827.099999999999
$ws.Range("N132").Value = -31713.332

$ws = $wb.Worksheets.Item("WVR")
# Row 31: Whatchoo Talking About / Cotton Doublet Vest of Crafting
$ws.Range("H31").Value = 1200
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 1731.3334
$ws.Range("I122").Value = 1779.5161
$ws.Range("J122").Value = 1595.5454
$ws.Range("K122").Value = 5338.5483
$ws.Range("L122").Value = 4786.6362
$ws.Range("M122").Value = -2888.5483
$ws.Range("N122").Value = -9686.636200000001

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 1481.48
$ws.Range("I132").Value = 899.35596
$ws.Range("J132").Value = 2319.1707
$ws.Range("K132").Value = 2698.06788
$ws.Range("L132").Value = 6957.5121
$ws.Range("M132").Value = -168.0678800000001
$ws.Range("N132").Value = -12017.5121
